# Edit: "Added links and addressed comments"
# Target shape: "TextBox 32" on slide 1 (the "new rows appended to a unbounded
# table" caption next to the structured-streaming diagram).
#
# 1) Resize/reposition the textbox (it grew wider once "to a unbounded table"
#    was split into "to " / "an " / "unbounded table" so the word wrap -- the
#    box uses spAutoFit -- settled on a slightly different box).
# 2) Fix the "a unbounded" -> "an unbounded" typo, and split the final run
#    into three runs ("to ", "an ", "unbounded table") so that "to " and
#    "an " can later carry their own hyperlinks independent of the rest of
#    the sentence.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "TextBox 32") {
        $sh = $cand
    }
}

# --- reposition / resize (values chosen so the single-precision round trip
# used by the host lands exactly on the target EMU numbers) ---
$sh.Left  = 546.320556640625
$sh.Width = 184.1812744140625

# --- fix the text & split it into three runs ---
$tr = $sh.TextFrame.TextRange
$lastPara = $tr.Paragraphs(5, 1)
$lastPara.Text = "unbounded table"
$null = $lastPara.InsertBefore("an ")
$null = $lastPara.InsertBefore("to ")
